$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 262, shifting rows 262:332 down to 263:333
$ws.Rows.Item(262).Insert()

# Populate the new row 262 with the new data entry
$ws.Cells.Item(262, 1).Value = 9
$ws.Cells.Item(262, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(262, 3).Value = "Metropolitana"
$ws.Cells.Item(262, 4).Value = 44855
$ws.Cells.Item(262, 5).Value = 13
$ws.Cells.Item(262, 6).Value = 300000001
$ws.Cells.Item(262, 7).Value = "Rabanito"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 5000
$ws.Cells.Item(262, 11).Value = 4000
$ws.Cells.Item(262, 12).Value = 4000
$ws.Cells.Item(262, 13).Value = 4000
$ws.Cells.Item(262, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(262, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(262, 16).Value = 40
$ws.Cells.Item(262, 17).Value = 100
$ws.Cells.Item(262, 18).Value = "Hortaliza"
